$d = $word.ActiveDocument
$d.Content.Find.Execute("Huella de Carbono2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Huella de Carbono223", 2)
